# "revised how to read excel"
#
# The sheet used to carry the event name ("30sスピード") and the
# division name ("3年生部門") on their own little row (row 3), sitting
# above the real header row (row 4: 順位/氏名/団体/回数). That makes the
# data awkward to read/import as a table, so those two labels are moved
# onto the header row itself, in new columns E ("30sスピード") and F
# ("3年生部門") - and then repeated down every data row so each record
# is fully self-describing on its own line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old standalone label row - its content now lives in E:F.
$ws.Rows(3).ClearContents()

# Re-add the "30sスピード" / "3年生部門" labels: once on the header row,
# and once per data row (rows 4 through 23).
for ($r = 4; $r -le 23; $r++) {
    $ws.Range("E$r").Value = "30sスピード"
    $ws.Range("F$r").Value = "3年生部門"
}

# Leave the selection where the editor ended up after making this change.
$ws.Range("I25").Select() | Out-Null
